{"js": "// Rename the template's bracketed placeholders to their final field names and\n// relocate the stray \"_GoBack\" bookmark to the very last paragraph (matching\n// the position Word leaves it in after the last edit of the document).\n\nconst body = context.document.body;\n\n// 1. Simple placeholder-text renames: [uniqueReference] -> [UniqueId], etc.\nconst renames = [\n  [\"uniqueReference\", \"UniqueId\"],\n  [\"manager\", \"ManagerName\"],\n  [\"startDate\", \"StartDate\"],\n  [\"endDate\", \"EndDate\"],\n  [\"daysRequestedCount\", \"DaysrequestedCount\"],\n  [\"additionalNotes\", \"AdditionalNotes\"],\n  [\"employee\", \"EmployeeName\"],\n];\n\nfor (const [oldText, newText] of renames) {\n  const results = body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n\n  if (results.items.length === 0) {\n    continue;\n  }\n  for (const r of results.items) {\n    r.insertText(newText, \"Replace\");\n  }\n  await context.sync();\n}\n\n// 2. Move the \"_GoBack\" bookmark: remove it from wherever it currently sits\n// (right after \"Additional Notes: [AdditionalNotes\") and re-insert it in the\n// last paragraph of the document (an empty paragraph at the very end).\nconst existing = context.document.getBookmarkRangeOrNullObject(\"_GoBack\");\nawait context.sync();\nif (!existing.isNullObject) {\n  context.document.deleteBookmark(\"_GoBack\");\n}\n\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items\");\nawait context.sync();\n\nconst lastParagraph = paragraphs.items[paragraphs.items.length - 1];\nlastParagraph.getRange().insertBookmark(\"_GoBack\");\nawait context.sync();\n", "ps1": "# Rename the template's bracketed placeholders to their final field names and\n# relocate the stray \"_GoBack\" bookmark to the very last paragraph (matching\n# the position Word leaves it in after the last edit of the document).\n\n$d = $word.ActiveDocument\n\n# 1. Simple placeholder-text renames: [uniqueReference] -> [UniqueId], etc.\n$renames = @(\n    @(\"uniqueReference\", \"UniqueId\"),\n    @(\"manager\", \"ManagerName\"),\n    @(\"startDate\", \"StartDate\"),\n    @(\"endDate\", \"EndDate\"),\n    @(\"daysRequestedCount\", \"DaysrequestedCount\"),\n    @(\"additionalNotes\", \"AdditionalNotes\"),\n    @(\"employee\", \"EmployeeName\")\n)\n\nforeach ($pair in $renames) {\n    $findText = $pair[0]\n    $replaceText = $pair[1]\n\n    $rng = $d.Content\n    $rng.Find.ClearFormatting()\n    $rng.Find.Replacement.ClearFormatting()\n    $rng.Find.Execute(\n        $findText,\n        $true,          # MatchCase\n        $false,         # MatchWholeWord\n        $false,         # MatchWildcards\n        $false,         # MatchSoundsLike\n        $false,         # MatchAllWordForms\n        $true,          # Forward\n        1,              # Wrap -> wdFindContinue\n        $false,         # Format\n        $replaceText,   # ReplaceWith\n        2               # Replace -> wdReplaceAll\n    )\n}\n\n# 2. Move the \"_GoBack\" bookmark: remove it from wherever it currently sits\n# (right after \"Additional Notes: [AdditionalNotes\") and re-insert it in the\n# last paragraph of the document (an empty paragraph at the very end).\nif ($d.Bookmarks.Exists(\"_GoBack\")) {\n    $d.Bookmarks.Item(\"_GoBack\").Delete()\n}\n\n$lastParagraph = $d.Paragraphs.Item($d.Paragraphs.Count)\n$d.Bookmarks.Add(\"_GoBack\", $lastParagraph.Range)\n"}
